# Completes draft ApplicationFlow; Adds user stories re: flagging content
#
# Appends a batch of new "Task" notes (column D) to the time-log sheet,
# picking the log back up a few rows below the prior content (with a few
# intentionally blank rows acting as paragraph breaks, matching the
# author's original layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "Mon Am 45 min starting app flow"
$ws.Range("D29").Value = "15 min dao video reference cave of coding"
$ws.Range("D30").Value = "11:50 - flow"

$ws.Range("D32").Value = "Need to"
$ws.Range("D33").Value = "Revise search interface and flow to match."
$ws.Range("D34").Value = "How will event be triggered to search for matching careers?"
$ws.Range("D35").Value = "How will event  be triggered to select matching career and get the median income?"
$ws.Range("D36").Value = "Let's skilp the income populating on the form and just show it in the results."

$ws.Range("D39").Value = "Paragraph or summary needs to show user what career they searched on (if they did), what the median income is, and"
$ws.Range("D40").Value = "all the other stuff already mentioned in diagram."

$ws.Range("D42").Value = "Need to be consistent across flow/user stories/screens about which capabilities are given to people with simple survey profiles"
$ws.Range("D43").Value = "and which are only for those who have written their story."

$ws.Range("D45").Value = "Maybe users have read-only access to everything until they have written their story?"
$ws.Range("D46").Value = "Or maybe they only see charts in search and have read-only access to forums?"

$ws.Range("D48").Value = "Need a way to go back to the thread list after viewing a specific thread."

$ws.Range("D50").Value = "Add sitewide search function for special topic, maybe on nav bar for signed in users."

# Row 42 wraps onto a second line at the column's width, same as the other
# two-line entries already in the sheet (ht="30").
$ws.Rows.Item(42).RowHeight = 30

# Scroll the sheet so the newly added block is in view, and leave the
# cursor on the first empty row right after it (mirrors the author's
# final cursor position after typing the last entry).
$ws.Range("D51").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 3
